$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B" = 8.220758201098199; "C" = 4.93467453886397; "D" = 4.881175997749213; "E" = 12.67062366166256; "F" = 23.90263805532801; "I" = 20.26921622773893; "K" = 8.14779931892344; "M" = 13.46619406058953; "N" = 18.33807785221151; "O" = 21.41567170606819 }
    3 = @{ "B" = 7.92806961936692; "C" = 4.73823452452827; "D" = 4.820547725536577; "E" = 12.4573648603666; "F" = 23.90998162060156; "I" = 20.341392761627; "K" = 7.951884181708333; "M" = 13.29288771035219; "N" = 18.3963928304861; "O" = 21.46902490471089 }
    4 = @{ "B" = 7.744015815610905; "C" = 4.612328853306642; "D" = 4.782338047141671; "E" = 12.32896841116271; "F" = 23.92074991205313; "I" = 20.38948115228745; "K" = 7.830426438292053; "M" = 13.18861510946581; "N" = 18.43385432712245; "O" = 21.50634686456419 }
    5 = @{ "B" = 7.668040314494464; "C" = 4.559733325199814; "D" = 4.76652783375876; "E" = 12.27736165416721; "F" = 23.92671142334098; "I" = 20.41002495245291; "K" = 7.78071174132229; "M" = 13.14671187391397; "N" = 18.44953781459846; "O" = 21.52270097668872 }
    6 = @{ "B" = 7.655369648302027; "C" = 4.550923444353755; "D" = 4.763888331842502; "E" = 12.26883781789243; "F" = 23.92779631768746; "I" = 20.41349342447848; "K" = 7.772445587380139; "M" = 13.13979089846945; "N" = 18.45216730573743; "O" = 21.52548565680572 }
    7 = @{ "B" = 7.742994952899007; "C" = 4.611624686104484; "D" = 4.782125784066163; "E" = 12.32826942670874; "F" = 23.9208239421829; "I" = 20.38975437859706; "K" = 7.82975675888469; "M" = 13.18804753758584; "N" = 18.43406414753473; "O" = 21.50656278846693 }
    8 = @{ "B" = 8.120808665561519; "C" = 4.868063808845049; "D" = 4.86048019692163; "E" = 12.59660952864927; "F" = 23.90387090462704; "I" = 20.29331926060177; "K" = 8.08053256372455; "M" = 13.40602517712256; "N" = 18.35784193395293; "O" = 21.43311955355933 }
    9 = @{ "B" = 8.822558160894069; "C" = 5.327376979641368; "D" = 5.00599396224014; "E" = 13.13965335086819; "F" = 23.92027909411717; "I" = 20.13418302714175; "K" = 8.560009463997368; "M" = 13.84815894264739; "N" = 18.22145116971608; "O" = 21.3253969143252 }
    10 = @{ "B" = 9.30875641816662; "C" = 5.636506771014067; "D" = 5.107512961311321; "E" = 13.54438947012217; "F" = 23.96253064923669; "I" = 20.03560008479012; "K" = 8.900969501184454; "M" = 14.17885303274385; "N" = 18.1291368826326; "O" = 21.26850238276326 }
    11 = @{ "B" = 9.522602549761755; "C" = 5.770699597202014; "D" = 5.152440654569366; "E" = 13.72883384716153; "F" = 23.98827233850564; "I" = 19.99474654074817; "K" = 9.052899603146649; "M" = 14.32991228314893; "N" = 18.08883679835117; "O" = 21.24747247171424 }
    12 = @{ "B" = 9.60245948572814; "C" = 5.820570788022919; "D" = 5.169266447463502; "E" = 13.79865038772116; "F" = 23.99895294670428; "I" = 19.97985158671269; "K" = 9.109920524694024; "M" = 14.38715058152361; "N" = 18.07381852327948; "O" = 21.24020805139585 }
    13 = @{ "B" = 9.585311760542202; "C" = 5.809872433777492; "D" = 5.165651153475751; "E" = 13.78361659136027; "F" = 23.99661129065864; "I" = 19.98303387008159; "K" = 9.097663640595323; "M" = 14.37482256752969; "N" = 18.0770422104146; "O" = 21.2417414583372 }
    14 = @{ "B" = 9.529195273929759; "C" = 5.774821566050394; "D" = 5.153828716785764; "E" = 13.7345786916434; "F" = 23.98913237620653; "I" = 19.99350958470743; "K" = 9.057601305006093; "M" = 14.33462090967067; "N" = 18.08759638232132; "O" = 21.2468608038288 }
    15 = @{ "B" = 9.494674315882326; "C" = 5.753228332655684; "D" = 5.146562530240079; "E" = 13.70453564594259; "F" = 23.98467263681512; "I" = 20.00000123439884; "K" = 9.032993762969795; "M" = 14.30999925344997; "N" = 18.09409266380947; "O" = 21.25008764069561 }
    16 = @{ "B" = 9.29462758615864; "C" = 5.627605760554498; "D" = 5.10455099138426; "E" = 13.53233566204107; "F" = 23.9609791935481; "I" = 20.0383504077557; "K" = 8.890971813416011; "M" = 14.16898885208928; "N" = 18.1318045821771; "O" = 21.2699744976465 }
    17 = @{ "B" = 9.169977243428658; "C" = 5.548877887368751; "D" = 5.078452030784241; "E" = 13.42672583092446; "F" = 23.94811108637703; "I" = 20.06289984857985; "K" = 8.802992225524294; "M" = 14.08260451854619; "N" = 18.15537273149127; "O" = 21.2834182067339 }
    18 = @{ "B" = 9.097595428052369; "C" = 5.502991265674557; "D" = 5.06332298787998; "E" = 13.36601708395036; "F" = 23.94132386531529; "I" = 20.07739573771577; "K" = 8.75209293475185; "M" = 14.03298054403693; "N" = 18.16908803538276; "O" = 21.29160726349697 }
    19 = @{ "B" = 9.072972548324882; "C" = 5.487351689385169; "D" = 5.058180573839156; "E" = 13.34547054718802; "F" = 23.93913144773722; "I" = 20.08236828355329; "K" = 8.734810267819695; "M" = 14.01619104507919; "N" = 18.17375923509661; "O" = 21.29445830826885 }
    20 = @{ "B" = 9.183318029396089; "C" = 5.557321313964576; "D" = 5.081242534019262; "E" = 13.43796508805351; "F" = 23.94941738623807; "I" = 20.06024762563121; "K" = 8.812388826863486; "M" = 14.09179425833044; "N" = 18.15284735940273; "O" = 21.28193983399473 }
    21 = @{ "B" = 9.545708997115664; "C" = 5.785142637366083; "D" = 5.157306392372437; "E" = 13.74898367551838; "F" = 23.99130384374569; "I" = 19.99041698698749; "K" = 9.06938289035255; "M" = 14.3464285771307; "N" = 18.0844897939358; "O" = 21.24533814259657 }
    22 = @{ "B" = 9.775987070629306; "C" = 5.928523438850735; "D" = 5.205923443704933; "E" = 13.9520528783209; "F" = 24.02411303692283; "I" = 19.94813307071274; "K" = 9.26438895883498; "M" = 14.51302903587995; "N" = 18.04122708891403; "O" = 21.22549252207365 }
    23 = @{ "B" = 9.653704424264181; "C" = 5.8525086545157; "D" = 5.180078061281628; "E" = 13.84371353768609; "F" = 24.00610685233294; "I" = 19.97039343474979; "K" = 9.147701070672316; "M" = 14.42411222347832; "N" = 18.06418829578033; "O" = 21.23571114922781 }
    24 = @{ "B" = 9.177288893293643; "C" = 5.5535059857793; "D" = 5.07998133445313; "E" = 13.43288378996904; "F" = 23.94882490502122; "I" = 20.06144550552796; "K" = 8.808141611681162; "M" = 14.08763944967563; "N" = 18.15398856399855; "O" = 21.28260677325956 }
    25 = @{ "B" = 8.637517193199793; "C" = 5.20798468406523; "D" = 4.967545162502847; "E" = 12.99141156243948; "F" = 23.9105281031351; "I" = 20.17401833771684; "K" = 8.432024242600122; "M" = 13.72729739277816; "N" = 18.25695665149558; "O" = 21.35063941834073 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

Write-Output "Updated $($data.Count) rows"